$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2314.2856
$ws.Range("I29").Value = 1961.2
$ws.Range("J29").Value = 3197
$ws.Range("K29").Value = 5883.6
$ws.Range("L29").Value = 9591
$ws.Range("M29").Value = -5602.6
$ws.Range("N29").Value = -10153
$ws.Range("H38").Value = 2873.6
$ws.Range("I38").Value = 194.5
$ws.Range("J38").Value = 4659.6665
$ws.Range("K38").Value = 583.5
$ws.Range("L38").Value = 13978.9995
$ws.Range("M38").Value = -211.5
$ws.Range("N38").Value = -14722.9995
$ws.Range("H116").Value = 3268.2
$ws.Range("I116").Value = 3268.2
$ws.Range("K116").Value = 3268.2
$ws.Range("M116").Value = 173.8000000000002
$ws.Range("H135").Value = 1473.1538
$ws.Range("I135").Value = 1646.6666
$ws.Range("J135").Value = 1082.75
$ws.Range("K135").Value = 14819.9994
$ws.Range("L135").Value = 9744.75
$ws.Range("M135").Value = -12284.9994
$ws.Range("N135").Value = -14814.75
$ws.Range("H137").Value = 7953.909
$ws.Range("I137").Value = 9999
$ws.Range("J137").Value = 4375
$ws.Range("K137").Value = 29997
$ws.Range("L137").Value = 13125
$ws.Range("M137").Value = -27447
$ws.Range("N137").Value = -18225

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2519.611
$ws.Range("J45").Value = 2842.4285
$ws.Range("L45").Value = 2842.4285
$ws.Range("N45").Value = -3596.4285
$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050
$ws.Range("H132").Value = 2902.5334
$ws.Range("I132").Value = 2322
$ws.Range("J132").Value = 4499
$ws.Range("K132").Value = 6966
$ws.Range("L132").Value = 13497
$ws.Range("M132").Value = -4436
$ws.Range("N132").Value = -18557

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1467.875
$ws.Range("I86").Value = 1463.2858
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 1463.2858
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -340.2858000000001
$ws.Range("N86").Value = -3746
$ws.Range("H89").Value = 1467.875
$ws.Range("I89").Value = 1463.2858
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 7316.429
$ws.Range("L89").Value = 7500
$ws.Range("M89").Value = -1700.429
$ws.Range("N89").Value = -18732

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4034.8865
$ws.Range("I31").Value = 2540.8147
$ws.Range("J31").Value = 6407.8237
$ws.Range("K31").Value = 2540.8147
$ws.Range("L31").Value = 6407.8237
$ws.Range("M31").Value = -2245.8147
$ws.Range("N31").Value = -6997.8237
$ws.Range("H34").Value = 4034.8865
$ws.Range("I34").Value = 2540.8147
$ws.Range("J34").Value = 6407.8237
$ws.Range("K34").Value = 2540.8147
$ws.Range("L34").Value = 6407.8237
$ws.Range("M34").Value = -2338.8147
$ws.Range("N34").Value = -6811.8237
$ws.Range("H58").Value = 2319.6667
$ws.Range("I58").Value = 2286
$ws.Range("K58").Value = 2286
$ws.Range("M58").Value = -2083
$ws.Range("H105").Value = 1597.4
$ws.Range("I105").Value = 1597.4
$ws.Range("K105").Value = 1597.4
$ws.Range("M105").Value = 149.5999999999999
$ws.Range("H132").Value = 4304.476
$ws.Range("I132").Value = 3863.2
$ws.Range("K132").Value = 11589.6
$ws.Range("M132").Value = -9059.599999999999
$ws.Range("H136").Value = 2319.6667
$ws.Range("I136").Value = 2286
$ws.Range("K136").Value = 6858
$ws.Range("M136").Value = -4308

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 829.6842
$ws.Range("J5").Value = 947.7
$ws.Range("L5").Value = 2843.1
$ws.Range("N5").Value = -3067.1
$ws.Range("H34").Value = 1978
$ws.Range("J34").Value = 2795
$ws.Range("L34").Value = 8385
$ws.Range("N34").Value = -8553
$ws.Range("H119").Value = 5001
$ws.Range("I119").Value = 5001
$ws.Range("K119").Value = 15003
$ws.Range("M119").Value = -10165
$ws.Range("H135").Value = 829.6842
$ws.Range("J135").Value = 947.7
$ws.Range("L135").Value = 8529.300000000001
$ws.Range("N135").Value = -13599.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 24398
$ws.Range("J96").Value = 11597
$ws.Range("L96").Value = 11597
$ws.Range("N96").Value = -17089
$ws.Range("H132").Value = 5402.2
$ws.Range("I132").Value = 5146.5713
$ws.Range("K132").Value = 15439.7139
$ws.Range("M132").Value = -12909.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3992.3076
$ws.Range("I22").Value = 3900
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 3900
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = -3605
$ws.Range("N22").Value = -4590
$ws.Range("H27").Value = 3992.3076
$ws.Range("I27").Value = 3900
$ws.Range("J27").Value = 4000
$ws.Range("K27").Value = 3900
$ws.Range("L27").Value = 4000
$ws.Range("M27").Value = -3793
$ws.Range("N27").Value = -4214
$ws.Range("H43").Value = 901599.75
$ws.Range("J43").Value = 901599.75
$ws.Range("L43").Value = 901599.75
$ws.Range("N43").Value = -901985.75
$ws.Range("H55").Value = 358.6
$ws.Range("I55").Value = 358.6
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 358.6
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -185.6
$ws.Range("N55").Value = ""
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = ""
$ws.Range("H132").Value = 5223.8335
$ws.Range("I132").Value = 4820.8
$ws.Range("K132").Value = 14462.4
$ws.Range("M132").Value = -11932.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 18334.666
$ws.Range("I6").Value = 15005
$ws.Range("J6").Value = 19999.5
$ws.Range("K6").Value = 15005
$ws.Range("L6").Value = 19999.5
$ws.Range("M6").Value = -14890
$ws.Range("N6").Value = -20229.5
$ws.Range("H81").Value = 5712.1665
$ws.Range("I81").Value = 5712.1665
$ws.Range("K81").Value = 11424.333
$ws.Range("M81").Value = -10363.333
$ws.Range("H84").Value = 5712.1665
$ws.Range("I84").Value = 5712.1665
$ws.Range("K84").Value = 57121.665
$ws.Range("M84").Value = -51817.665
$ws.Range("H122").Value = 996
$ws.Range("I122").Value = 996
$ws.Range("K122").Value = 2988
$ws.Range("M122").Value = -538
$ws.Range("H132").Value = 3353.48
$ws.Range("I132").Value = 3271.75
$ws.Range("J132").Value = 3498.7778
$ws.Range("K132").Value = 9815.25
$ws.Range("L132").Value = 10496.3334
$ws.Range("M132").Value = -7285.25
$ws.Range("N132").Value = -15556.3334
$ws.Range("H136").Value = 12286
$ws.Range("I136").Value = 12485.923
$ws.Range("K136").Value = 37457.769
$ws.Range("M136").Value = -34907.769
